$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.083.43"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "1.909.38"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -1.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.88"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -1.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4831"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3822"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07355"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9345"
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.80"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07802"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "1.900.91"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.509"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.627"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.42"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  -1.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008827"
$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "28.097.80"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.147"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "2.150.93"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.98"
$ws.Range("E25").Value = "  +2.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.925"
$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.105"
$ws.Range("E28").Value = "  +4.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.32"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.954"
$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08915"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.337"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("E33").Value = "  +2.76%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.687"
$ws.Range("E34").Value = "  +1.39%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7664"
$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.608"
$ws.Range("E36").Value = "  -2.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02049"
$ws.Range("E37").Value = "  -0.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.102"
$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05297"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5487"
$ws.Range("E40").Value = "  +2.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.980"
$ws.Range("E41").Value = "  -0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.021"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1522"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.454"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4831"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.22"
$ws.Range("E47").Value = "  +3.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.38"
$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06108"
$ws.Range("E51").Value = "  +0.16%  "
